# Updated cryptos list with GitHub Actions
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the crypto
# table, and swaps the EnergySwap/Aptos rows (47/48) per the latest data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.318.01'
$ws.Range("D3").Value = '1.870.55'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.91'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4706'
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2885'
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06621'
$ws.Range("E9").Value = '  +1.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.79'
$ws.Range("E10").Value = '  +1.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08028'
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.49'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '1.871.10'
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.147'
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6867'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '271.27'
$ws.Range("E16").Value = '  -2.49%  '
$ws.Range("D17").Value = '30.309.70'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.14'
$ws.Range("E18").Value = '  +3.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007727'
$ws.Range("E19").Value = '  +5.82%  '
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '2.116.15'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.307'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.224'
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.428'
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.76'
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.96'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.958'
$ws.Range("E28").Value = '  +1.79%  '
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09898'
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.378'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.468'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.084'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04710'
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7025'
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.649'
$ws.Range("E39").Value = '  +2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.307'
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.90'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.961'
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8437'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4176'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.28'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.087'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.206'
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '927.15'
$ws.Range("E49").Value = '  -5.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.49'
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05683'
$ws.Range("E51").Value = '  +0.62%  '
